$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A91").Value = 0.21
$ws.Range("B91").Value = 0.41
$ws.Range("C91").Value = 0.21
$ws.Range("D91").Value = 0.799999
$ws.Range("E91").Value = 0.5999989999999999
$ws.Range("F91").Value = 0.799999
$ws.Range("G91").Value = 0.4
$ws.Range("H91").Value = 0.799999
$ws.Range("I91").Value = 9.2725017399974
$ws.Range("J91").Value = "query"
